$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.208.69"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.03%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.849.11"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.33%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.83%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.7016"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.41%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07717"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3070"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.19%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.65"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.55%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07826"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.24%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "92.93"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.74%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.136"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.82%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.839.47"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.09%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6868"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.02%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.611"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.40%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008313"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.23%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.201.32"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.11%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "242.17"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.85%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.088.99"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.20%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.75"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.44%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9998"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.03%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.520"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.17%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.001"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.08%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1512"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.77%  "

$ws.Range("E26").Value = "  -0.84%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.837"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.09%  "

$ws.Range("E28").Value = "  -0.86%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.542"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.99%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.230"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.10%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.185"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.20%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.205"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.95%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05122"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.55%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7921"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.20%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.921"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.32%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.147"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.55%  "

$ws.Range("E37").Value = "  -0.49%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.325.04"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.88%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01872"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.57%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.715"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.20%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9553"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.47%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.074"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +9.80%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "107.42"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.14%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.15%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "9.709"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.84%  "

$ws.Range("E46").Value = "  -0.97%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.988.95"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.18%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5183"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.06%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "64.30"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.42%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.764"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.83%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.997"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.17%  "
